$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trade row (row 8) appended below the existing HZNP noun-trade history.
$ws.Range("A8").Value = 8760
$ws.Range("B8").Value = 8880.7800000000007
$ws.Range("C8").Value = 19.170000000000002
$ws.Range("D8").Value = 19.43
$ws.Range("E8").Value = $true
$ws.Range("F8").Value = 1.36
$ws.Range("G8").Value = 42609.488726851851
$ws.Range("G8").NumberFormat = "m/d/yy h:mm"
$ws.Range("H8").Value = $false
